$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.146.04'
$ws.Range('E2').Value = '  +1.44%  '
$ws.Range('D3').Value = '2.589.15'
$ws.Range('E3').Value = '  +2.98%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '315.55'
$c.Style = "Normal"
$ws.Range('E5').Value = '  -0.70%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '97.50'
$c.Style = "Normal"
$ws.Range('E6').Value = '  +3.29%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('E9').Value = '  +1.31%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '35.62'
$c.Style = "Normal"
$ws.Range('E10').Value = '  -0.14%  '
$ws.Range('E11').Value = '  +0.12%  '
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '7.50'
$c.Style = "Normal"
$ws.Range('E12').Value = '  -0.36%  '
$ws.Range('D13').Value = '2.987.63'
$ws.Range('E13').Value = '  +2.98%  '
$ws.Range('E14').Value = '  -0.67%  '
$ws.Range('D15').Value = '2.517.59'
$ws.Range('E15').Value = '  +1.06%  '
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '15.25'
$c.Style = "Normal"
$ws.Range('E16').Value = '  +0.31%  '
$ws.Range('E17').Value = '  +0.10%  '
$ws.Range('D18').Value = '43.220.60'
$ws.Range('E18').Value = '  +1.44%  '
$ws.Range('E19').Value = '  +2.50%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '12.68'
$c.Style = "Normal"
$ws.Range('E20').Value = '  -1.66%  '
$ws.Range('E21').Value = '  +0.96%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '69.56'
$c.Style = "Normal"
$ws.Range('E22').Value = '  +0.43%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '251.37'
$c.Style = "Normal"
$ws.Range('E23').Value = '  +0.19%  '
$ws.Range('E24').Value = '  +0.18%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '2.08'
$c.Style = "Normal"
$ws.Range('E25').Value = '  +3.14%  '
$ws.Range('E26').Value = '  +2.57%  '
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '2.43'
$c.Style = "Normal"
$ws.Range('E28').Value = '  -0.81%  '
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '40.87'
$c.Style = "Normal"
$ws.Range('E29').Value = '  -0.53%  '
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '10.29'
$c.Style = "Normal"
$ws.Range('E30').Value = '  +0.75%  '
$ws.Range('E31').Value = '  -1.03%  '
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '156.84'
$c.Style = "Normal"
$ws.Range('E32').Value = '  -0.17%  '
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '3.44'
$c.Style = "Normal"
$ws.Range('E33').Value = '  +5.55%  '
$ws.Range('E34').Value = '  +1.99%  '
$ws.Range('E35').Value = '  +3.57%  '
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '2.70'
$c.Style = "Normal"
$ws.Range('E36').Value = '  +2.92%  '
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '18.81'
$c.Style = "Normal"
$ws.Range('E37').Value = '  -1.54%  '
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '2.51'
$c.Style = "Normal"
$ws.Range('E38').Value = '  +9.83%  '
$ws.Range('E39').Value = '  +1.97%  '
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '0.118'
$c.Style = "Normal"
$ws.Range('E40').Value = '  +0.47%  '
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '23.19'
$c.Style = "Normal"
$ws.Range('E41').Value = '  -1.51%  '
$ws.Range('E42').Value = '  +5.46%  '
$ws.Range('E43').Value = '  +1.01%  '
$ws.Range('E44').Value = '  -0.19%  '
$ws.Range('D45').Value = '2.013.57'
$ws.Range('E45').Value = '  +0.00%  '
$ws.Range('E46').Value = '  -2.09%  '
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '9.00'
$c.Style = "Normal"
$ws.Range('E47').Value = '  +1.09%  '
$ws.Range('D48').Value = '2.839.52'
$ws.Range('E48').Value = '  +3.01%  '
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '83.26'
$c.Style = "Normal"
$ws.Range('E49').Value = '  -2.06%  '
$ws.Range('E50').Value = '  +4.44%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '74.40'
$c.Style = "Normal"
$ws.Range('E51').Value = '  -0.23%  '

Write-Host "Applied 77 cell updates"
